$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("G3 y G4")
$cos = $ws.ChartObjects()
$co = $cos.Item(2)
$chart = $co.Chart
$s2 = $chart.SeriesCollection(2)
try {
  Write-Host $s2.Interior
} catch {
  Write-Host "no interior: $_"
}
